$d = $word.ActiveDocument

# Use wdReplaceOne (1) rather than wdReplaceAll (2): several replacement
# strings end with the exact search string (e.g. " | Oct 2021" becomes
# "Senior Software Engineer | Oct 2021"), and a "replace all" pass keeps
# rescanning forward, re-matching the tail of text it just inserted and
# duplicating the prefix. Every target below (save one, handled specially)
# is a single, unique match in the document, so "replace first occurrence"
# is exactly what we want.
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 1) | Out-Null
}

# 1. Subtitle
Replace-Text "Fullstack .Net Developer" "Senior Software Engineer"

# 2. Betsson heading + title/date line
Replace-Text "24 Betsson" "Betsson Group"
Replace-Text " | Oct 2021" "Senior Software Engineer | Oct 2021"

# 3. Plejmo heading
Replace-Text "17 Plejmo" "Film2Home/Plejmo"

# " | Nov 2014" occurs twice (Plejmo's own entry and the later Film2Home AB
# entry), and the replacement text for the Plejmo one re-introduces the
# literal search string as its own suffix. Locate both occurrences by plain
# text offset first, then fix up the later (Film2Home AB) one via a
# range-restricted Find so it can't be confused with the other, then do the
# earlier (Plejmo) one - which is then the only remaining match.
$t = $d.Content.Text
$idxFirst = $t.IndexOf(" | Nov 2014")
$idxSecond = $t.IndexOf(" | Nov 2014", $idxFirst + 1)

$rSecond = $d.Range($idxSecond, $d.Content.End)
$rSecond.Find.Execute(" | Nov 2014", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "Fullstack Developer | Nov 2014", 1) | Out-Null

Replace-Text " | Nov 2014" "Fullstack Developer, DevOps | Nov 2014"

# 4. Adlibris heading + title/date line
Replace-Text "23 Adlibris" "Adlibris"
Replace-Text " | Feb 2021" "Fullstack Developer | Feb 2021"

# 5. Open Payments heading + title/date line
Replace-Text "21 Open Payments" "Open Payments"
Replace-Text " | Dec 2019" "Software Engineer | Dec 2019"

# 6. Henrik Becker Consulting AB
Replace-Text " | Jul 2017" "IT Consultant | Jul 2017"

# 7. Betsson Group AB
Replace-Text " | Dec 2023" "Senior Software Engineer | Dec 2023"

# 8. Magine TV AB
Replace-Text " | Sep 2015" "Fullstack Developer | Sep 2015"

# 9. Qbranch Stockholm AB (first instance, Apr 2014)
Replace-Text " | Apr 2014" "IT Consultant | Apr 2014"

# 10. Wasa Kredit AB
Replace-Text " | Mar 2012" "Lead Developer | Mar 2012"

# 11. Avega Group AB
Replace-Text " | Jan 2008" "IT Consultant | Jan 2008"

# 12. Aftonbladet Tillväxtteknik 2 AB
Replace-Text " | Dec 2007" "Fullstack Developer | Dec 2007"

# 13. Inverso International AB
Replace-Text " | Apr 2007" "IT Consultant | Apr 2007"

# 14. Qbranch Stockholm AB (second instance, Apr 2005)
Replace-Text " | Apr 2005" "IT Consultant | Apr 2005"

# 15. eWork
Replace-Text " | Jan 2002" "IT Consultant | Jan 2002"

# 16. IconMedialab AB
Replace-Text " | Sep 1999" "IT Consultant | Sep 1999"

# 17. Innitek AB
Replace-Text " | Oct 1998" "IT Consultant | Oct 1998"

# 18. Lidingö Stad
Replace-Text " | Jul 1997" "PC Technician | Jul 1997"

# 19. Konsumentföreningen Stockholm
Replace-Text " | Dec 1986" "Shop Assistant | Dec 1986"

# 20. Sollentuna Jazz Workshop
Replace-Text " | Jan 1994" "Study Circle Leader | Jan 1994"

# 21. Kantarellen Livs AB
Replace-Text " | Jun 1992" "Shop Assistant | Jun 1992"

# 22. Karena Zoo, Lek & Hobby
Replace-Text " | Aug 1986" "Shop Assistant | Aug 1986"

Write-Host "All replacements applied"
